$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format to prevent Excel auto-converting
# numeric-looking strings (e.g. "0.0410", "565.50") into numbers,
# which would silently drop significant trailing/represented digits.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '65.493.50'
$ws.Range("E2").Value = '  -3.35%  '
$ws.Range("D3").Value = '3.468.91'
$ws.Range("E3").Value = '  -2.16%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '565.50'
$ws.Range("E5").Value = '  +1.82%  '
$ws.Range("D6").Value = '175.69'
$ws.Range("E6").Value = '  -8.61%  '
$ws.Range("D7").Value = '0.628'
$ws.Range("E7").Value = '  +2.83%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E9").Value = '  -3.10%  '
$ws.Range("E10").Value = '  +1.10%  '
$ws.Range("D11").Value = '53.04'
$ws.Range("E11").Value = '  -6.32%  '
$ws.Range("D12").Value = '0.0000267'
$ws.Range("E12").Value = '  -2.54%  '
$ws.Range("E13").Value = '  -4.75%  '
$ws.Range("D14").Value = '4.021.04'
$ws.Range("E14").Value = '  -2.07%  '
$ws.Range("D15").Value = '3.464.63'
$ws.Range("E15").Value = '  -2.23%  '
$ws.Range("E16").Value = '  -0.60%  '
$ws.Range("D17").Value = '18.07'
$ws.Range("E17").Value = '  -2.00%  '
$ws.Range("D18").Value = '65.451.85'
$ws.Range("E18").Value = '  -3.50%  '
$ws.Range("D19").Value = '11.91'
$ws.Range("E19").Value = '  -0.44%  '
$ws.Range("D20").Value = '0.985'
$ws.Range("E20").Value = '  -1.67%  '
$ws.Range("D21").Value = '409.93'
$ws.Range("E21").Value = '  +0.36%  '
$ws.Range("D22").Value = '4.09'
$ws.Range("E22").Value = '  +2.81%  '
$ws.Range("D23").Value = '84.48'
$ws.Range("E23").Value = '  -0.94%  '
$ws.Range("D24").Value = '4.14'
$ws.Range("E24").Value = '  -3.53%  '
$ws.Range("D25").Value = '12.61'
$ws.Range("E25").Value = '  +5.01%  '
$ws.Range("D26").Value = '10.82'
$ws.Range("E26").Value = '  -5.53%  '
$ws.Range("E27").Value = '  -6.32%  '
$ws.Range("D28").Value = '8.85'
$ws.Range("E28").Value = '  +0.18%  '
$ws.Range("D29").Value = '29.87'
$ws.Range("E29").Value = '  -2.47%  '
$ws.Range("D30").Value = '611.69'
$ws.Range("E30").Value = '  -12.63%  '
$ws.Range("E31").Value = '  -9.15%  '
$ws.Range("D32").Value = '11.50'
$ws.Range("E32").Value = '  -2.78%  '
$ws.Range("E33").Value = '  -3.83%  '
$ws.Range("D34").Value = '58.66'
$ws.Range("E34").Value = '  -3.44%  '
$ws.Range("D35").Value = '0.149'
$ws.Range("E35").Value = '  +7.02%  '
$ws.Range("D36").Value = '0.999'
$ws.Range("E36").Value = '  -0.13%  '
$ws.Range("B37").Value = 'PEPE'
$ws.Range("C37").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D37").Value = '0.0₃0781'
$ws.Range("E37").Value = '  -6.36%  '
$ws.Range("B38").Value = 'Maker'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D38").Value = '3.353.94'
$ws.Range("E38").Value = '  +9.72%  '
$ws.Range("D39").Value = '36.43'
$ws.Range("E39").Value = '  -7.32%  '
$ws.Range("E40").Value = '  -6.13%  '
$ws.Range("D41").Value = '3.40'
$ws.Range("E41").Value = '  +1.04%  '
$ws.Range("D42").Value = '0.998'
$ws.Range("E42").Value = '  -0.18%  '
$ws.Range("D43").Value = '2.87'
$ws.Range("E43").Value = '  -4.51%  '
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").Value = '0.0410'
$ws.Range("E44").Value = '  -3.34%  '
$ws.Range("B45").Value = 'ApeXProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D45").Value = '3.20'
$ws.Range("E45").Value = '  -2.66%  '
$ws.Range("D46").Value = '2.47'
$ws.Range("E46").Value = '  -8.37%  '
$ws.Range("D47").Value = '2.67'
$ws.Range("E47").Value = '  -2.18%  '
$ws.Range("E48").Value = '  -0.34%  '
$ws.Range("D49").Value = '137.40'
$ws.Range("E49").Value = '  -2.30%  '
$ws.Range("D50").Value = '8.33'
$ws.Range("E50").Value = '  -9.09%  '
$ws.Range("E51").Value = '  +6.46%  '

# Restore default (unstyled) cell style on column D now that the
# text values are committed, matching the original workbook where
# these cells carry no explicit style.
$ws.Range("D2:D51").Style = "Normal"
